$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '36.881.37'
$ws.Range('E2').Value = '  -1.60%  '
$ws.Range('D3').Value = '2.019.72'
$ws.Range('E3').Value = '  -2.45%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '225.42'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.96%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.606'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.64%  '
$ws.Range('E7').Value = '  +0.08%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '54.58'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -4.78%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.378'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.74%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0786'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.99%  '
$ws.Range('E11').Value = '  -3.59%  '
$ws.Range('D12').Value = '2.320.01'
$ws.Range('E12').Value = '  -2.35%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '14.24'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -3.95%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '20.23'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.99%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.743'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.81%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.13'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -3.50%  '
$ws.Range('D17').Value = '2.012.07'
$ws.Range('E17').Value = '  -2.96%  '
$ws.Range('D18').Value = '36.832.87'
$ws.Range('E18').Value = '  -1.41%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.18'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +3.66%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '68.59'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.60%  '
$ws.Range('D21').Value = '0.0₃0818'
$ws.Range('E21').Value = '  -1.15%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '226.68'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.65%  '
$ws.Range('E23').Value = '  -0.07%  '
$ws.Range('E24').Value = '  +2.47%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.20'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -7.01%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '164.92'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.98%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.17'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -4.60%  '
$ws.Range('E28').Value = '  -6.03%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '18.66'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -4.06%  '
$ws.Range('E30').Value = '  -4.19%  '
$ws.Range('E31').Value = '  -4.98%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.45'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.55%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0614'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -3.04%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.42'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -4.34%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.34'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -5.11%  '
$ws.Range('E36').Value = '  +0.90%  '
$ws.Range('E37').Value = '  +0.06%  '
$ws.Range('E38').Value = '  -4.99%  '
$ws.Range('E39').Value = '  +0.26%  '
$ws.Range('B40').Value = 'Maker'
$ws.Range('C40').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D40').Value = '1.488.31'
$ws.Range('E40').Value = '  +1.39%  '
$ws.Range('B41').Value = 'InjectiveProtocol'
$ws.Range('C41').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '17.07'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.87%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0217'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -5.66%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0926'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.55%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '94.70'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -4.98%  '
$ws.Range('E45').Value = '  -5.51%  '
$ws.Range('E46').Value = '  -5.72%  '
$ws.Range('E47').Value = '  +1.45%  '
$ws.Range('E48').Value = '  -3.95%  '
$ws.Range('E49').Value = '  -0.93%  '
$ws.Range('D50').Value = '2.207.88'
$ws.Range('E50').Value = '  -2.27%  '
$ws.Range('E51').Value = '  -8.63%  '
